$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.414.34"
$ws.Range("E2").Value = "  +0.45%  "

# Row 3
$ws.Range("D3").Value = "3.543.12"
$ws.Range("E3").Value = "  +0.04%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "600.74"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.38%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.43%  "

# Row 7
$ws.Range("D7").Value = "3.543.62"
$ws.Range("E7").Value = "  +0.08%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.09%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.498"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.87%  "

# Row 10
$ws.Range("E10").Value = "  +0.35%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.16"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.43%  "

# Row 12
$ws.Range("E12").Value = "  +0.50%  "

# Row 13
$ws.Range("D13").Value = "4.148.26"
$ws.Range("E13").Value = "  +0.12%  "

# Row 14
$ws.Range("E14").Value = "  +1.13%  "

# Row 15
$ws.Range("E15").Value = "  +0.31%  "

# Row 16
$ws.Range("D16").Value = "3.548.72"
$ws.Range("E16").Value = "  +0.38%  "

# Row 17
$ws.Range("E17").Value = "  -0.06%  "

# Row 18
$ws.Range("D18").Value = "65.488.72"
$ws.Range("E18").Value = "  +0.58%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.81"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.99%  "

# Row 20
$ws.Range("E20").Value = "  +2.45%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.76"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.92%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "393.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.78%  "

# Row 23
$ws.Range("E23").Value = "  +1.44%  "

# Row 24
$ws.Range("D24").Value = "3.689.64"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "75.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.53%  "

# Row 26
$ws.Range("E26").Value = "  +0.01%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000115"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.57%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.94"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.56%  "

# Row 29
$ws.Range("E29").Value = "  +15.71%  "

# Row 30
$ws.Range("E30").Value = "  -2.26%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.38%  "

# Row 32
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.51"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.85%  "

# Row 33
$ws.Range("D33").Value = "3.550.92"
$ws.Range("E33").Value = "  -0.24%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "24.24"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.56%  "

# Row 35
$ws.Range("E35").Value = "  +0.02%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.148"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.31%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.35"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +6.46%  "

# Row 38
$ws.Range("E38").Value = "  +0.41%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "169.71"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.56%  "

# Row 40
$ws.Range("E40").Value = "  -0.37%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0831"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.35%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.829"
$ws.Range("D42").Style = "Normal"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "26.04"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.17%  "

# Row 44
$ws.Range("E44").Value = "  +4.21%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "42.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.97%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.07%  "

# Row 47
$ws.Range("E47").Value = "  +0.73%  "

# Row 48
$ws.Range("E48").Value = "  +0.05%  "

# Row 49
$ws.Range("E49").Value = "  +1.30%  "

# Row 50
$ws.Range("D50").Value = "2.381.54"
$ws.Range("E50").Value = "  -0.45%  "

# Row 51
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0270"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.26%  "
